# Femacal de La Calera - Ajo: add a new weekly price observation.
# A new row is inserted at row 520 (pushing the existing rows 520:551
# down to 521:552); the new row carries the latest weekly reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(520).Insert()

$ws.Cells.Item(520, 1).Value  = 3
$ws.Cells.Item(520, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(520, 3).Value  = "Coquimbo"
$ws.Cells.Item(520, 4).Value  = 44826
$ws.Cells.Item(520, 5).Value  = 5
$ws.Cells.Item(520, 6).Value  = 100112003
$ws.Cells.Item(520, 7).Value  = "Ajo"
$ws.Cells.Item(520, 8).Value  = "Chino"
$ws.Cells.Item(520, 9).Value  = "Primera"
$ws.Cells.Item(520, 10).Value = 85
$ws.Cells.Item(520, 11).Value = 16000
$ws.Cells.Item(520, 12).Value = 16500
$ws.Cells.Item(520, 13).Value = 16235
$ws.Cells.Item(520, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(520, 15).Value = "China"
$ws.Cells.Item(520, 16).Value = 1624
$ws.Cells.Item(520, 17).Value = 10
$ws.Cells.Item(520, 18).Value = "Hortaliza"
